$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the absolute path recorded for this workbook (x15ac:absPath)
# ---------------------------------------------------------------------------
try {
    $wb.Path = "F:\My Drive\MudrikLab020818\Experiments_new\subliminal_priming_w_motion_capture\development\"
} catch {
}

# ---------------------------------------------------------------------------
# 2) Restructure the columns:
#    target_x/y/z/timecourse (AA:AD) and prime_x/y/z/timecourse (AI:AL, i.e.
#    the block right after target_rt) each get split into a "_to" and a
#    "_from" variant (8 columns instead of 4), while the pas_x/y/z/timecourse
#    block is removed entirely.
# ---------------------------------------------------------------------------

# Insert 4 new columns right after target_timecourse (AD) to make room for
# the target "_from" columns.
$ws.Range("AE1:AH1").EntireColumn.Insert()

# Insert 4 new columns right after prime_timecourse (now AP) to make room
# for the prime "_from" columns.
$ws.Range("AQ1:AT1").EntireColumn.Insert()

# Delete the (now shifted) pas_x/pas_y/pas_z/pas_timecourse columns (AX:BA).
$ws.Range("AX1:BA1").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 3) Rewrite the header (row 1) and description (row 2) text for the
#    restructured trajectory columns.
# ---------------------------------------------------------------------------

# target trajectory "to screen"
$ws.Range("AA1").Value = "target_x_to"
$ws.Range("AB1").Value = "target_y_to"
$ws.Range("AC1").Value = "target_z_to"
$ws.Range("AD1").Value = "target_timecourse_to"
$ws.Range("AA2").Value = "X trajectory to screen"
$ws.Range("AB2").Value = "Y trajectory to screen"
$ws.Range("AC2").Value = "Z trajectory to screen"
$ws.Range("AD2").Value = "time of each trajectory sample (sec) when reaching to screen"

# target trajectory "from screen to start point"
$ws.Range("AE1").Value = "target_x_from"
$ws.Range("AF1").Value = "target_y_from"
$ws.Range("AG1").Value = "target_z_from"
$ws.Range("AH1").Value = "target_timecourse_from"
$ws.Range("AE2").Value = "X trajectory from screen to start point"
$ws.Range("AF2").Value = "Y trajectory from screen to start point"
$ws.Range("AG2").Value = "Z trajectory from screen to start point"
$ws.Range("AH2").Value = "time of each trajectory sample (sec) when returning from screen"

# prime trajectory "to screen"
$ws.Range("AM1").Value = "prime_x_to"
$ws.Range("AN1").Value = "prime_y_to"
$ws.Range("AO1").Value = "prime_z_to"
$ws.Range("AP1").Value = "prime_timecourse_to"
$ws.Range("AM2").Value = "X trajectory to screen"
$ws.Range("AN2").Value = "Y trajectory to screen"
$ws.Range("AO2").Value = "Z trajectory to screen"
$ws.Range("AP2").Value = "time of each trajectory sample (sec) when reaching to screen"

# prime trajectory "from screen to start point"
$ws.Range("AQ1").Value = "prime_x_from"
$ws.Range("AR1").Value = "prime_y_from"
$ws.Range("AS1").Value = "prime_z_from"
$ws.Range("AT1").Value = "prime_timecourse_from"
$ws.Range("AQ2").Value = "X trajectory from screen to start point"
$ws.Range("AR2").Value = "Y trajectory from screen to start point"
$ws.Range("AS2").Value = "Z trajectory from screen to start point"
$ws.Range("AT2").Value = "time of each trajectory sample (sec) when returning from screen"

# ---------------------------------------------------------------------------
# 4) Misc formatting / view updates
# ---------------------------------------------------------------------------

# Header row is now taller.
$ws.Rows.Item(1).RowHeight = 28.5

# Update the view: scroll / selection.
[void]$ws.Range("W17").Select()
try {
    $excel.ActiveWindow.ScrollColumn = 15
    $excel.ActiveWindow.ScrollRow = 1
} catch {
}
